$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: odds refresh (Q4, R4) ---
$ws.Range("Q4").Value = 1.8
$ws.Range("R4").Value = 2

# --- Row 5: odds refresh ---
$ws.Range("G5").Value = 2.1
$ws.Range("H5").Value = 3.3
$ws.Range("I5").Value = 3.4
$ws.Range("J5").Value = 2.75
$ws.Range("K5").Value = 2.2
$ws.Range("L5").Value = 4
$ws.Range("Q5").Value = 1.98
$ws.Range("R5").Value = 1.88
$ws.Range("X5").Value = 10
$ws.Range("Z5").Value = 19
$ws.Range("AA5").Value = 17
$ws.Range("AB5").Value = 26
$ws.Range("AC5").Value = 10
$ws.Range("AD5").Value = 6.5
$ws.Range("AF5").Value = 41
$ws.Range("AH5").Value = 11
$ws.Range("AK5").Value = 41
$ws.Range("AL5").Value = 29
$ws.Range("AN5").Value = 4
$ws.Range("AO5").Value = 11
$ws.Range("AW5").Value = 5.5
$ws.Range("BD5").Value = 126

# --- Insert two new rows at position 6, pushing old rows 6-8 down to 8-10 ---
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).Insert()

# --- New Row 6: xripwnPs / Plymouth - Portsmouth ---
$ws.Range("B6:C6").NumberFormat = "@"
$ws.Range("A6").Value = "xripwnPs"
$ws.Range("B6").Value = "05/11/2024"
$ws.Range("C6").Value = "17:00"
$ws.Range("D6").Value = "ENGLAND - CHAMPIONSHIP"
$ws.Range("E6").Value = "Plymouth"
$ws.Range("F6").Value = "Portsmouth"
$ws.Range("G6").Value = 2.35
$ws.Range("H6").Value = 3.5
$ws.Range("I6").Value = 2.8
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 2.25
$ws.Range("L6").Value = 3.4
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 13
$ws.Range("O6").Value = 1.22
$ws.Range("P6").Value = 4.33
$ws.Range("Q6").Value = 1.73
$ws.Range("R6").Value = 2.1
$ws.Range("S6").Value = 1.33
$ws.Range("T6").Value = 3.25
$ws.Range("U6").Value = 1.62
$ws.Range("V6").Value = 2.2
$ws.Range("W6").Value = 10
$ws.Range("X6").Value = 13
$ws.Range("Y6").Value = 9.5
$ws.Range("Z6").Value = 23
$ws.Range("AA6").Value = 17
$ws.Range("AB6").Value = 23
$ws.Range("AC6").Value = 13
$ws.Range("AD6").Value = 7
$ws.Range("AE6").Value = 12
$ws.Range("AF6").Value = 41
$ws.Range("AG6").Value = 151
$ws.Range("AH6").Value = 11
$ws.Range("AI6").Value = 15
$ws.Range("AJ6").Value = 11
$ws.Range("AK6").Value = 29
$ws.Range("AL6").Value = 21
$ws.Range("AM6").Value = 26
$ws.Range("AN6").Value = 4.5
$ws.Range("AO6").Value = 13
$ws.Range("AP6").Value = 21
$ws.Range("AQ6").Value = 41
$ws.Range("AR6").Value = 51
$ws.Range("AS6").Value = 151
$ws.Range("AT6").Value = 2.75
$ws.Range("AU6").Value = 7.5
$ws.Range("AV6").Value = 41
$ws.Range("AW6").Value = 5
$ws.Range("AX6").Value = 15
$ws.Range("AY6").Value = 21
$ws.Range("AZ6").Value = 51
$ws.Range("BA6").Value = 51
$ws.Range("BB6").Value = 126
$ws.Range("BC6").Value = 451
$ws.Range("BD6").Value = 151

# --- New Row 7: KhRHJuMt / Swansea - Watford ---
$ws.Range("B7:C7").NumberFormat = "@"
$ws.Range("A7").Value = "KhRHJuMt"
$ws.Range("B7").Value = "05/11/2024"
$ws.Range("C7").Value = "17:00"
$ws.Range("D7").Value = "ENGLAND - CHAMPIONSHIP"
$ws.Range("E7").Value = "Swansea"
$ws.Range("F7").Value = "Watford"
$ws.Range("G7").Value = 1.85
$ws.Range("H7").Value = 3.6
$ws.Range("I7").Value = 4
$ws.Range("J7").Value = 2.5
$ws.Range("K7").Value = 2.3
$ws.Range("L7").Value = 4
$ws.Range("M7").Value = 1.04
$ws.Range("N7").Value = 13
$ws.Range("O7").Value = 1.22
$ws.Range("P7").Value = 4.33
$ws.Range("Q7").Value = 1.67
$ws.Range("R7").Value = 2.2
$ws.Range("S7").Value = 1.33
$ws.Range("T7").Value = 3.25
$ws.Range("U7").Value = 1.62
$ws.Range("V7").Value = 2.2
$ws.Range("W7").Value = 9
$ws.Range("X7").Value = 10
$ws.Range("Y7").Value = 8.5
$ws.Range("Z7").Value = 17
$ws.Range("AA7").Value = 15
$ws.Range("AB7").Value = 21
$ws.Range("AC7").Value = 13
$ws.Range("AD7").Value = 7
$ws.Range("AE7").Value = 13
$ws.Range("AF7").Value = 41
$ws.Range("AG7").Value = 151
$ws.Range("AH7").Value = 13
$ws.Range("AI7").Value = 21
$ws.Range("AJ7").Value = 13
$ws.Range("AK7").Value = 41
$ws.Range("AL7").Value = 29
$ws.Range("AM7").Value = 34
$ws.Range("AN7").Value = 4
$ws.Range("AO7").Value = 10
$ws.Range("AP7").Value = 19
$ws.Range("AQ7").Value = 34
$ws.Range("AR7").Value = 51
$ws.Range("AS7").Value = 101
$ws.Range("AT7").Value = 3.25
$ws.Range("AU7").Value = 7.5
$ws.Range("AV7").Value = 41
$ws.Range("AW7").Value = 6
$ws.Range("AX7").Value = 21
$ws.Range("AY7").Value = 26
$ws.Range("AZ7").Value = 67
$ws.Range("BA7").Value = 81
$ws.Range("BB7").Value = 151
$ws.Range("BC7").Value = 501
$ws.Range("BD7").Value = 151
